{"js": "// Update the lattice-multiplication exercise table: every cell's\n// multiplication problem / grid is replaced with a new one.  The\n// table shape (5 rows x 3 cols) is unchanged - only the text inside\n// each cell is different.  Each cell holds one run made of 5 lines\n// (joined by <w:br/>) which Office.js exposes through `text`/\n// `insertText` using \"\\v\" (vertical tab) as the line separator.\nconst newValues = [\n  [\n    \"46 x 59\\v  5    9\\v  ----\\v4|    |\\v6|    |\",\n    \"29 x 27\\v  2    7\\v  ----\\v2|    |\\v9|    |\",\n    \"51 x 70\\v  7    0\\v  ----\\v5|    |\\v1|    |\",\n  ],\n  [\n    \"16 x 15\\v  1    5\\v  ----\\v1|    |\\v6|    |\",\n    \"85 x 19\\v  1    9\\v  ----\\v8|    |\\v5|    |\",\n    \"51 x 60\\v  6    0\\v  ----\\v5|    |\\v1|    |\",\n  ],\n  [\n    \"99 x 86\\v  8    6\\v  ----\\v9|    |\\v9|    |\",\n    \"19 x 59\\v  5    9\\v  ----\\v1|    |\\v9|    |\",\n    \"26 x 16\\v  1    6\\v  ----\\v2|    |\\v6|    |\",\n  ],\n  [\n    \"84 x 41\\v  4    1\\v  ----\\v8|    |\\v4|    |\",\n    \"76 x 80\\v  8    0\\v  ----\\v7|    |\\v6|    |\",\n    \"11 x 45\\v  4    5\\v  ----\\v1|    |\\v1|    |\",\n  ],\n  [\n    \"65 x 38\\v  3    8\\v  ----\\v6|    |\\v5|    |\",\n    \"51 x 72\\v  7    2\\v  ----\\v5|    |\\v1|    |\",\n    \"59 x 71\\v  7    1\\v  ----\\v5|    |\\v9|    |\",\n  ],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the lattice-multiplication exercise table: every cell's\n# multiplication problem / grid is replaced with a new one. The table\n# shape (5 rows x 3 cols) does not change - only the text inside each\n# cell differs. Each cell is a single paragraph/run whose 5 lines are\n# joined by line breaks (<w:br/>), which the Word COM model represents\n# as Chr(11) (vertical tab) inside Range.Text.\n$vt = [char]11\n\n$newValues = @(\n    @(\n        @(\"46 x 59\", \"  5    9\", \"  ----\", \"4|    |\", \"6|    |\"),\n        @(\"29 x 27\", \"  2    7\", \"  ----\", \"2|    |\", \"9|    |\"),\n        @(\"51 x 70\", \"  7    0\", \"  ----\", \"5|    |\", \"1|    |\")\n    ),\n    @(\n        @(\"16 x 15\", \"  1    5\", \"  ----\", \"1|    |\", \"6|    |\"),\n        @(\"85 x 19\", \"  1    9\", \"  ----\", \"8|    |\", \"5|    |\"),\n        @(\"51 x 60\", \"  6    0\", \"  ----\", \"5|    |\", \"1|    |\")\n    ),\n    @(\n        @(\"99 x 86\", \"  8    6\", \"  ----\", \"9|    |\", \"9|    |\"),\n        @(\"19 x 59\", \"  5    9\", \"  ----\", \"1|    |\", \"9|    |\"),\n        @(\"26 x 16\", \"  1    6\", \"  ----\", \"2|    |\", \"6|    |\")\n    ),\n    @(\n        @(\"84 x 41\", \"  4    1\", \"  ----\", \"8|    |\", \"4|    |\"),\n        @(\"76 x 80\", \"  8    0\", \"  ----\", \"7|    |\", \"6|    |\"),\n        @(\"11 x 45\", \"  4    5\", \"  ----\", \"1|    |\", \"1|    |\")\n    ),\n    @(\n        @(\"65 x 38\", \"  3    8\", \"  ----\", \"6|    |\", \"5|    |\"),\n        @(\"51 x 72\", \"  7    2\", \"  ----\", \"5|    |\", \"1|    |\"),\n        @(\"59 x 71\", \"  7    1\", \"  ----\", \"5|    |\", \"9|    |\")\n    )\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $lines = $rowValues[$c - 1]\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = [string]::Join($vt, $lines)\n    }\n}\n"}
